$d = $word.ActiveDocument

# 1. The "\cite{Julia-by-example}" text had accidentally been split across
#    three runs, with a bogus proofErr gramStart/gramEnd pair wrapped around
#    the literal "cite{" piece. Re-join it into a single run of plain text
#    (this also drops the now-unneeded proofErr markers, since Find/Replace
#    rewrites the matched range as one run).
$old = "familiarity with the language\cite{Julia-by-example}. Key areas like File I/O and Random Number generation were then learned as these are fundamentally important for the AST and file generation. The complexities of creating the abstract portions of an AST in Julia were then attempted."
$new = $old
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
Write-Output ("Merged citation run: " + $found)

# 2. Append two blank paragraphs after the final paragraph, matching the
#    formatting (ListBullet style, no numbering, hanging indent, minorHAnsi
#    cstheme font) already used for the other empty paragraphs in this
#    document.
$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter()
$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter()
Write-Output ("Paragraph count now: " + $d.Paragraphs.Count)
